$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Starting Mode" -> "Base Mode"  (only the word "Starting" changes)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Starting Mode", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Base Mode", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Locate the paragraph that carries the "_GoBack" bookmark (it is an
#    empty ilvl=3 list item right before "Laser Contest w/ STAR ...").
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bmPos = $bm.Start

$bookmarkParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $bmPos -and $p.Range.End -gt $bmPos) {
        $bookmarkParaIndex = $i
        break
    }
}

$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

# Insert a fresh blank ilvl=3 list item directly above it (this becomes the
# new standalone empty bullet that precedes "Laser Contest...").
$bookmarkPara.Range.InsertParagraphBefore() | Out-Null

# The original (bookmarked) empty paragraph was pushed down by one; find it
# again via the bookmark, which has not moved.
$bm = $d.Bookmarks.Item("_GoBack")
$bmPos = $bm.Start

$bookmarkParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $bmPos -and $p.Range.End -gt $bmPos) {
        $bookmarkParaIndex = $i
        break
    }
}
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

# Remove the bookmark, then delete this now-redundant empty paragraph
# entirely (its content merges away, leaving the "Laser Contest..."
# paragraph immediately following the new blank bullet).
$bm.Delete()

$mergeStart = $bookmarkPara.Range.Start
$mergeEnd = $bookmarkPara.Range.End
$mergeRange = $d.Range($mergeStart, $mergeEnd)
$mergeRange.Delete()

# The paragraph that now begins at $mergeStart is "Laser Contest w/ STAR
# (Accuracy Challenge)"; re-add the "_GoBack" bookmark at its very start.
$laserParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $mergeStart) {
        $laserParaIndex = $i
        break
    }
}
$laserPara = $d.Paragraphs.Item($laserParaIndex)
$newBookmarkRange = $d.Range($laserPara.Range.Start, $laserPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange) | Out-Null
